$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the numeric-looking Price/Volume columns so Excel
# does not auto-convert values like "1.001" or "30.207.55" into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.207.55'
$ws.Range("D3").Value = '1.859.39'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '235.47'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.4665'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '0.2838'
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("D9").Value = '0.06508'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '21.45'
$ws.Range("E10").Value = '  +6.75%  '
$ws.Range("D11").Value = '0.07901'
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '96.99'
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '1.858.20'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '5.146'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '0.6761'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("D16").Value = '278.04'
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '30.212.80'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '13.60'
$ws.Range("E18").Value = '  +8.00%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '5.376'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.104.65'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").Value = '0.000007289'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '6.126'
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '166.41'
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").Value = '9.170'
$ws.Range("E26").Value = '  -1.52%  '
$ws.Range("D27").Value = '19.00'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = '1.921'
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("D29").Value = '1.384'
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("D30").Value = '0.09703'
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("D31").Value = '4.361'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '1.475'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").Value = '4.022'
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").Value = '0.04709'
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("D35").Value = '1.125'
$ws.Range("E35").Value = '  +2.56%  '
$ws.Range("D36").Value = '0.7038'
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("D37").Value = '2.709'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '0.01852'
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = '2.620'
$ws.Range("E39").Value = '  +4.32%  '
$ws.Range("D40").Value = '6.308'
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").Value = '74.35'
$ws.Range("E41").Value = '  +3.26%  '
$ws.Range("D42").Value = '1.946'
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '0.8492'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4157'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '103.23'
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("D47").Value = '985.38'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.138'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.279'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("D50").Value = '33.93'
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05645'
$ws.Range("E51").Value = '  +0.13%  '

# Remove the temporary text-number-format so the cells keep the default (unstyled)
# appearance that the original workbook used for these cells.
$ws.Range("D2:E51").ClearFormats()
